$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 24
$ws.Range("D2").Value = "A las (hora) de (mes) del año (año), el suscrito Ing. Mauricio Arroyo Herrera, Director de la Escuela de Ingeniería en Computación en atención al caso de modificación de acta del estudiante ESTUDIANTE 2, carné 2015567890, sobre el curso IC5821 REQUERIMIENTOS DE SOFTWARE, grupo 1, del (número) semestre del año (año), resuelvo:"
$ws.Range("E2").Value = "Por un error involuntario, no se tramitó la inclusión del estudiante ESTUDIANTE 2 con identificación 2015567890 en el curso IC5821 REQUERIMIENTOS DE SOFTWARE, grupo 1 impartido por el profesor ALICIA SALAZAR HERNÁNDEZ en el (semestre) semestre del año (año) y al final del semestre, el profesor no pudo registrar en el acta la calificación obtenida por el estudiante."
$ws.Range("F2").Value = "Autorizar la modificación del acta del curso IC5821 REQUERIMIENTOS DE SOFTWARE, grupo 1 en el (semestre) semestre del año (año) impartido por el profesor ALICIA SALAZAR HERNÁNDEZ para incluir al estudiante ESTUDIANTE 2, carné 2015567890 con una nota de aprobación de (nota por escrito) (# de nota)."
$ws.Range("H2").Value = "Después de haber realizado la investigación del caso, y consultado al profesor ALICIA SALAZAR HERNÁNDEZ, quien impartió el curso IC5821 REQUERIMIENTOS DE SOFTWARE, grupo 1 en el (semestre) semestre del año (año) se logra comprobar que el estudiante ESTUDIANTE 2, carné 2015567890, efectivamente aprobó el curso con una nota de (nota por escrito) (# de nota), por lo que esta Dirección solicita gestionar la modificación del acta correspondiente."
